$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

$wsALC.Range("H19").Value = 7182797.5
$wsALC.Range("I19").Value = 5870135.5
$wsALC.Range("J19").Value = 9092124
$wsALC.Range("K19").Value = 5870135.5
$wsALC.Range("L19").Value = 9092124
$wsALC.Range("M19").Value = -5869960.5
$wsALC.Range("N19").Value = -9092474
$wsALC.Range("H28").Value = 467.04
$wsALC.Range("I28").Value = 330.21054
$wsALC.Range("K28").Value = 330.21054
$wsALC.Range("M28").Value = 154.78946
$wsALC.Range("H92").Value = 2021.2858
$wsALC.Range("I92").Value = 438.5
$wsALC.Range("J92").Value = 4131.6665
$wsALC.Range("K92").Value = 438.5
$wsALC.Range("L92").Value = 4131.6665
$wsALC.Range("M92").Value = 809.5
$wsALC.Range("N92").Value = -6627.6665
$wsALC.Range("H98").Value = 921.1539
$wsALC.Range("I98").Value = 847.5
$wsALC.Range("K98").Value = 847.5
$wsALC.Range("M98").Value = 650.5
$wsALC.Range("H116").Value = 5131698.5
$wsALC.Range("I116").Value = 6413781.5
$wsALC.Range("J116").Value = 3366.3333
$wsALC.Range("K116").Value = 6413781.5
$wsALC.Range("L116").Value = 3366.3333
$wsALC.Range("M116").Value = -6410339.5
$wsALC.Range("N116").Value = -10250.3333
$wsALC.Range("H122").Value = 921.1539
$wsALC.Range("I122").Value = 847.5
$wsALC.Range("K122").Value = 2542.5
$wsALC.Range("M122").Value = -92.5
$wsALC.Range("H135").Value = 500.69565
$wsALC.Range("I135").Value = 262
$wsALC.Range("J135").Value = 1177
$wsALC.Range("K135").Value = 2358
$wsALC.Range("L135").Value = 10593
$wsALC.Range("M135").Value = 177
$wsALC.Range("N135").Value = -15663
$wsARM.Range("H61").Value = 3639.1333
$wsARM.Range("I61").Value = 2718.6
$wsARM.Range("J61").Value = 5480.2
$wsARM.Range("K61").Value = 2718.6
$wsARM.Range("L61").Value = 5480.2
$wsARM.Range("M61").Value = -2506.6
$wsARM.Range("N61").Value = -5904.2
$wsARM.Range("H132").Value = 2614.5806
$wsARM.Range("I132").Value = 2601.4285
$wsARM.Range("J132").Value = 2642.2
$wsARM.Range("K132").Value = 7804.2855
$wsARM.Range("L132").Value = 7926.599999999999
$wsARM.Range("M132").Value = -5274.2855
$wsARM.Range("N132").Value = -12986.6
$wsARM.Range("H136").Value = 3639.1333
$wsARM.Range("I136").Value = 2718.6
$wsARM.Range("J136").Value = 5480.2
$wsARM.Range("K136").Value = 8155.799999999999
$wsARM.Range("L136").Value = 16440.6
$wsARM.Range("M136").Value = -5605.799999999999
$wsARM.Range("N136").Value = -21540.6
$wsARM.Range("H139").Value = 50142.332
$wsARM.Range("J139").Value = 50142.332
$wsARM.Range("L139").Value = 50142.332
$wsARM.Range("N139").Value = -60422.332
$wsBSM.Range("H94").Value = 875.6286
$wsBSM.Range("I94").Value = 729.5517
$wsBSM.Range("J94").Value = 1581.6666
$wsBSM.Range("K94").Value = 729.5517
$wsBSM.Range("L94").Value = 1581.6666
$wsBSM.Range("M94").Value = -278.5517
$wsBSM.Range("N94").Value = -2483.6666
$wsBSM.Range("H99").Value = 1867.3549
$wsBSM.Range("I99").Value = 1303.45
$wsBSM.Range("K99").Value = 1303.45
$wsBSM.Range("M99").Value = 194.55
$wsBSM.Range("H105").Value = 2106
$wsBSM.Range("I105").Value = 1728.9524
$wsBSM.Range("J105").Value = 2985.7778
$wsBSM.Range("K105").Value = 1728.9524
$wsBSM.Range("L105").Value = 2985.7778
$wsBSM.Range("M105").Value = 18.0476000000001
$wsBSM.Range("N105").Value = -6479.7778
$wsCRP.Range("H31").Value = 5050.1787
$wsCRP.Range("I31").Value = 4417.1763
$wsCRP.Range("J31").Value = 6028.4546
$wsCRP.Range("K31").Value = 4417.1763
$wsCRP.Range("L31").Value = 6028.4546
$wsCRP.Range("M31").Value = -4122.1763
$wsCRP.Range("N31").Value = -6618.4546
$wsCRP.Range("H34").Value = 5050.1787
$wsCRP.Range("I34").Value = 4417.1763
$wsCRP.Range("J34").Value = 6028.4546
$wsCRP.Range("K34").Value = 4417.1763
$wsCRP.Range("L34").Value = 6028.4546
$wsCRP.Range("M34").Value = -4215.1763
$wsCRP.Range("N34").Value = -6432.4546
$wsCRP.Range("H99").Value = 3125.0588
$wsCRP.Range("I99").Value = 2759.3333
$wsCRP.Range("J99").Value = 4002.8
$wsCRP.Range("K99").Value = 2759.3333
$wsCRP.Range("L99").Value = 4002.8
$wsCRP.Range("M99").Value = -1261.3333
$wsCRP.Range("N99").Value = -6998.8
$wsCRP.Range("H126").Value = 3125.0588
$wsCRP.Range("I126").Value = 2759.3333
$wsCRP.Range("J126").Value = 4002.8
$wsCRP.Range("K126").Value = 8277.999899999999
$wsCRP.Range("L126").Value = 12008.4
$wsCRP.Range("M126").Value = -5807.999899999999
$wsCRP.Range("N126").Value = -16948.4
$wsCRP.Range("H132").Value = 1763.7878
$wsCRP.Range("I132").Value = 1381.6957
$wsCRP.Range("K132").Value = 4145.0871
$wsCRP.Range("M132").Value = -1615.0871
$wsCRP.Range("H134").Value = 3930.2273
$wsCRP.Range("I134").Value = 2185.5
$wsCRP.Range("K134").Value = 6556.5
$wsCRP.Range("M134").Value = -4021.5
$wsCUL.Range("H113").Value = 1567933.4
$wsCUL.Range("I113").Value = 3448741
$wsCUL.Range("J113").Value = 593.6667
$wsCUL.Range("K113").Value = 10346223
$wsCUL.Range("L113").Value = 1781.0001
$wsCUL.Range("M113").Value = -10344053
$wsCUL.Range("N113").Value = -6121.0001
$wsCUL.Range("H131").Value = 708.1404
$wsCUL.Range("I131").Value = 307.89474
$wsCUL.Range("J131").Value = 908.2632
$wsCUL.Range("K131").Value = 923.6842200000001
$wsCUL.Range("L131").Value = 2724.7896
$wsCUL.Range("M131").Value = 4116.31578
$wsCUL.Range("N131").Value = -12804.7896
$wsLTW.Range("H22").Value = 435459.88
$wsLTW.Range("J22").Value = 1034.4
$wsLTW.Range("L22").Value = 1034.4
$wsLTW.Range("N22").Value = -1624.4
$wsLTW.Range("H27").Value = 435459.88
$wsLTW.Range("J27").Value = 1034.4
$wsLTW.Range("L27").Value = 1034.4
$wsLTW.Range("N27").Value = -1248.4
$wsLTW.Range("H46").Value = 795.35
$wsLTW.Range("I46").Value = 691.7273
$wsLTW.Range("J46").Value = 922
$wsLTW.Range("K46").Value = 691.7273
$wsLTW.Range("L46").Value = 922
$wsLTW.Range("M46").Value = -503.7273
$wsLTW.Range("N46").Value = -1298
$wsLTW.Range("H125").Value = 49905
$wsLTW.Range("J125").Value = 49905
$wsLTW.Range("L125").Value = 49905
$wsLTW.Range("N125").Value = -59745
$wsLTW.Range("H136").Value = 1900.4445
$wsLTW.Range("I136").Value = 1280.8
$wsLTW.Range("J136").Value = 2675
$wsLTW.Range("K136").Value = 3842.4
$wsLTW.Range("L136").Value = 8025
$wsLTW.Range("M136").Value = -1292.4
$wsLTW.Range("N136").Value = -13125
$wsWVR.Range("H138").Value = 38718.332
$wsWVR.Range("J138").Value = 38718.332
$wsWVR.Range("L138").Value = 38718.332
$wsWVR.Range("N138").Value = -48998.332
